$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2 updates
$ws.Range("I2").Value = 0.05501054111033029
$ws.Range("K2").Value = 169.56

# Row 3 updates
$ws.Range("C3").Value = 0.8846153846153846
$ws.Range("H3").Value = 0.8016011138183083
$ws.Range("I3").Value = 0.06219772500067548
$ws.Range("J3").Value = 0.7692307692307693
$ws.Range("K3").Value = 190.5769230769231
$ws.Range("Q3").Value = 32
$ws.Range("R3").Value = 40
$ws.Range("S3").Value = 90
$ws.Range("T3").Value = 155
$ws.Range("U3").Value = 247
$ws.Range("V3").Value = 2815
$ws.Range("W3").Value = 2807
$ws.Range("X3").Value = 2757
$ws.Range("Y3").Value = 2692
$ws.Range("Z3").Value = 2600
$ws.Range("AF3").Value = 0.98876
$ws.Range("AG3").Value = 0.98595
$ws.Range("AH3").Value = 0.968388
$ws.Range("AI3").Value = 0.945557
$ws.Range("AJ3").Value = 0.913242
